# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45170 (2023-09-01) to serial date 45174 (2023-09-05)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
